$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date column C, rows 2-11, is updated
# from serial 45221 (2023-10-22) to serial 45224 (2023-10-25).
$ws.Range("C2:C11").Value2 = 45224
